# Apply updated cryptocurrency market data (price + 1h volume change).
# Column D values are forced to Text (leading apostrophe, like typing
# '1.00 into Excel) so formatted numeric-looking strings such as "1.00",
# "0.220" or thousand-dot-separated prices like "70.873.46" keep their
# exact original text instead of being coerced into a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "70.873.46"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "'" + "3.635.56"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("D4").Value = "'" + "0.999"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'" + "604.49"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'" + "199.74"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'" + "0.999"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'" + "0.220"
$ws.Range("E9").Value = "  +9.85%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "'" + "53.85"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "'" + "0.0000307"
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("D13").Value = "'" + "9.57"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "'" + "4.209.12"
$ws.Range("E14").Value = "  +3.57%  "
$ws.Range("D15").Value = "'" + "640.42"
$ws.Range("E15").Value = "  +8.04%  "
$ws.Range("D16").Value = "'" + "12.99"
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("D17").Value = "'" + "70.909.18"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "'" + "3.652.35"
$ws.Range("E18").Value = "  +4.19%  "
$ws.Range("D19").Value = "'" + "19.06"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "'" + "1.00"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").Value = "'" + "18.51"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "'" + "104.20"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "'" + "3.01"
$ws.Range("E26").Value = "  -4.48%  "
$ws.Range("D27").Value = "'" + "10.48"
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("D28").Value = "'" + "9.75"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("D29").Value = "'" + "33.97"
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("D30").Value = "'" + "4.79"
$ws.Range("E30").Value = "  +11.91%  "
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("D32").Value = "'" + "12.26"
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").Value = "'" + "63.42"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "'" + "4.017.33"
$ws.Range("E35").Value = "  +8.32%  "
$ws.Range("D36").Value = "'" + "0.0₃0880"
$ws.Range("E36").Value = "  +6.36%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'" + "3.04"
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'" + "511.99"
$ws.Range("E39").Value = "  +7.78%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'" + "0.390"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'" + "36.71"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "'" + "3.55"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").Value = "'" + "0.0462"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").Value = "'" + "3.52"
$ws.Range("E45").Value = "  +7.10%  "
$ws.Range("D46").Value = "'" + "2.98"
$ws.Range("E46").Value = "  +6.28%  "
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("D48").Value = "'" + "8.66"
$ws.Range("E48").Value = "  +3.01%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'" + "2.96"
$ws.Range("E51").Value = "  +5.42%  "
